# TrafficSimulator/3Dobjects/3d_object_library.xlsx edit
# "Documentation. GUI development and new textures and positions."
#
# Data changes:
#   - M8            : "0.0" -> "20.0"   (new shared string)
#   - K24:K35       : "0.0" -> "55.0"   (new shared string, 12 rows)
# View changes:
#   - selection moves from E6 to N38 (and scrolled so column H is leftmost,
#     best-effort - scroll position persistence isn't reproducible here)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New position/texture values - writing literal text so they land as shared
# strings exactly like the authored workbook (t="s"), not as numbers.
$ws.Range("M8").Value = "20.0"
$ws.Range("K24:K35").Value = "55.0"

# Update the on-screen selection to match the saved view state.
$ws.Range("N38").Select()
